$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update time[us] value for row 6 (post_0_1_resize_hwc)
$ws.Range("S6").Value = 1664

# Update name column (B) values for rows 8-26 to reflect new ONNX tensor names
$ws.Range("B8").Value  = "input.4"
$ws.Range("B9").Value  = "input.8"
$ws.Range("B10").Value = "input.16"
$ws.Range("B11").Value = "input.20"
$ws.Range("B12").Value = "input.28"
$ws.Range("B13").Value = "input.32"
$ws.Range("B14").Value = "input.40"
$ws.Range("B15").Value = "input.44"
$ws.Range("B16").Value = "input.52"
$ws.Range("B17").Value = "input.56"
$ws.Range("B18").Value = "input.64"
$ws.Range("B19").Value = "input.68"
$ws.Range("B20").Value = "input.76"
$ws.Range("B21").Value = "input.88"
$ws.Range("B22").Value = "input.100"
$ws.Range("B24").Value = "input.116"
$ws.Range("B25").Value = "onnx::Concat_140"
$ws.Range("B26").Value = "input.132"

# Update Total time[us] value for row 32
$ws.Range("S32").Value = 28819
